$wb = $excel.ActiveWorkbook

# --- StrategyDictionaries sheet: update MACD_BB_Freeman JSON settings (ADX_THRESHOLD 30 -> 0) ---
$ws3 = $wb.Worksheets.Item("StrategyDictionaries")
$json = '{"MA_TYPE": "SMA", "MACD_FAST": 2, "MACD_SLOW": 11, "BB_PERIODS": 40, "BB_MULT": 2, "ADX": 3, "ADX_THRESHOLD": 0}'
$ws3.Range("B3").Value = $json

# --- Sheet1: update test case data ---
$ws1 = $wb.Worksheets.Item("Sheet1")

# Update the "From" date of test #1
$ws1.Range("D2").Value = 44470

# Fill the Optional Strategy Settings column with the (updated) MACD_BB_Freeman dictionary entry
$ws1.Range("K2").Value = $json

# Remove the now-redundant second test case row entirely (shifts later rows up by one)
$ws1.Rows.Item(3).Delete()

# --- Restore selections to match final state: StrategyDictionaries!B3 then back to Sheet1!D2 ---
$ws3.Activate()
$ws3.Range("B3").Select() | Out-Null

$ws1.Activate()
$ws1.Range("D2").Select() | Out-Null
